$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round row 5 values to the "custom accuracy" (2 decimal places) target values.
$ws.Range("B5").Value = 13.93
$ws.Range("C5").Value = 10.13
$ws.Range("D5").Value = 1.02
$ws.Range("E5").Value = 30.27
$ws.Range("F5").Value = 24.48
$ws.Range("G5").Value = 10.94
$ws.Range("H5").Value = 40.66
$ws.Range("I5").Value = 16.87
$ws.Range("J5").Value = 7.4
$ws.Range("K5").Value = 10.9
$ws.Range("L5").Value = 12.14
$ws.Range("M5").Value = 12.77
$ws.Range("N5").Value = 3.5
$ws.Range("O5").Value = 10.9
$ws.Range("P5").Value = 15.43
$ws.Range("Q5").Value = 9.31
$ws.Range("R5").Value = 0.77
$ws.Range("S5").Value = 0.65
$ws.Range("T5").Value = 158.61
$ws.Range("U5").Value = 30.4
$ws.Range("V5").Value = 10.06
$ws.Range("W5").Value = 20.29
$ws.Range("X5").Value = 10.64
$ws.Range("Y5").Value = 1.77
$ws.Range("Z5").Value = 19.81
$ws.Range("AA5").Value = 8.890000000000001
$ws.Range("AB5").Value = 7.93
$ws.Range("AC5").Value = 9.34
$ws.Range("AD5").Value = 12.75
$ws.Range("AE5").Value = 0.55
$ws.Range("AF5").Value = 36.77
$ws.Range("AG5").Value = 5.6
$ws.Range("AH5").Value = 12.58

# Remove row 6 (it no longer exists in the target data - 1000 rows worth of new data
# replaced the previous 5-row sample, shrinking the used range back to A1:AH5).
$ws.Rows(6).Delete()
